$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (15 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K28").Value = 4770.087
$ws.Range("N28").Value = -1608.2857
$ws.Range("I28").Value = 4770.087
$ws.Range("L28").Value = 638.2857
$ws.Range("J28").Value = 638.2857
$ws.Range("M28").Value = -4285.087
$ws.Range("H28").Value = 3806
$ws.Range("I125").Value = 127999
$ws.Range("H125").Value = 67387.44500000001
$ws.Range("K125").Value = 1151991
$ws.Range("M125").Value = -1149531
$ws.Range("M132").Value = -3879.736699999999
$ws.Range("I132").Value = 2136.5789
$ws.Range("K132").Value = 6409.736699999999
$ws.Range("H132").Value = 2443.84

# --- Sheet: ARM (47 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M2").Value = -138.90909
$ws.Range("L2").Value = 4230.5
$ws.Range("I2").Value = 251.90909
$ws.Range("H2").Value = 2146.476
$ws.Range("J2").Value = 4230.5
$ws.Range("K2").Value = 251.90909
$ws.Range("N2").Value = -4456.5
$ws.Range("H45").Value = 1616.6666
$ws.Range("K45").Value = 1425
$ws.Range("I45").Value = 1425
$ws.Range("M45").Value = -1048
$ws.Range("N63").Value = -9505.333500000001
$ws.Range("K63").Value = 1804.5
$ws.Range("H63").Value = 5601.8
$ws.Range("L63").Value = 8133.3335
$ws.Range("J63").Value = 8133.3335
$ws.Range("M63").Value = -1118.5
$ws.Range("I63").Value = 1804.5
$ws.Range("J66").Value = 8133.3335
$ws.Range("K66").Value = 9022.5
$ws.Range("H66").Value = 5601.8
$ws.Range("N66").Value = -47530.6675
$ws.Range("L66").Value = 40666.6675
$ws.Range("M66").Value = -5590.5
$ws.Range("I66").Value = 1804.5
$ws.Range("H80").Value = 54000
$ws.Range("N80").Value = -55996
$ws.Range("L80").Value = 54000
$ws.Range("J80").Value = 54000
$ws.Range("L83").Value = 162000
$ws.Range("J83").Value = 54000
$ws.Range("H83").Value = 54000
$ws.Range("N83").Value = -171984
$ws.Range("M116").Value = 2042.09091
$ws.Range("K116").Value = 251.90909
$ws.Range("J116").Value = 4230.5
$ws.Range("N116").Value = -8818.5
$ws.Range("L116").Value = 4230.5
$ws.Range("I116").Value = 251.90909
$ws.Range("H116").Value = 2146.476
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("H122").Value = 3041.889
$ws.Range("K122").Value = 9226.200000000001
$ws.Range("M122").Value = -6776.200000000001
$ws.Range("I122").Value = 3075.4
$ws.Range("J122").Value = 3000

# --- Sheet: BSM (56 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2146.476
$ws.Range("J3").Value = 4230.5
$ws.Range("K3").Value = 251.90909
$ws.Range("I3").Value = 251.90909
$ws.Range("N3").Value = -4458.5
$ws.Range("M3").Value = -137.90909
$ws.Range("L3").Value = 4230.5
$ws.Range("H62").Value = 51250
$ws.Range("L62").Value = 51250
$ws.Range("J62").Value = 51250
$ws.Range("N62").Value = -52622
$ws.Range("N65").Value = -160614
$ws.Range("H65").Value = 51250
$ws.Range("L65").Value = 153750
$ws.Range("J65").Value = 51250
$ws.Range("H82").Value = 54000
$ws.Range("I82").Value = 51000
$ws.Range("K82").Value = 51000
$ws.Range("M82").Value = -50617
$ws.Range("I85").Value = 51000
$ws.Range("H85").Value = 54000
$ws.Range("K85").Value = 51000
$ws.Range("M85").Value = -49674
$ws.Range("J86").Value = 62501676
$ws.Range("I86").Value = 1856.3334
$ws.Range("K86").Value = 1856.3334
$ws.Range("L86").Value = 62501676
$ws.Range("N86").Value = -62503922
$ws.Range("M86").Value = -733.3334
$ws.Range("H86").Value = 17243186
$ws.Range("K89").Value = 9281.666999999999
$ws.Range("J89").Value = 62501676
$ws.Range("M89").Value = -3665.666999999999
$ws.Range("N89").Value = -312519612
$ws.Range("L89").Value = 312508380
$ws.Range("I89").Value = 1856.3334
$ws.Range("H89").Value = 17243186
$ws.Range("J94").Value = 3998.8
$ws.Range("N94").Value = -4900.8
$ws.Range("H94").Value = 4227.5454
$ws.Range("L94").Value = 3998.8
$ws.Range("I94").Value = 4418.1665
$ws.Range("K94").Value = 4418.1665
$ws.Range("M94").Value = -3967.1665
$ws.Range("H99").Value = 4406.5264
$ws.Range("I99").Value = 3274.4
$ws.Range("M99").Value = -1776.4
$ws.Range("K99").Value = 3274.4
$ws.Range("L108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("H134").Value = 2396.7727
$ws.Range("I134").Value = 2301.5264
$ws.Range("K134").Value = 6904.5792
$ws.Range("M134").Value = -4369.5792

# --- Sheet: CRP (19 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J22").Value = 2160.625
$ws.Range("N22").Value = -2860.625
$ws.Range("K22").Value = 984.9286
$ws.Range("M22").Value = -634.9286
$ws.Range("L22").Value = 2160.625
$ws.Range("I22").Value = 984.9286
$ws.Range("H22").Value = 1412.4546
$ws.Range("I31").Value = 2101.6428
$ws.Range("H31").Value = 5116.8423
$ws.Range("M31").Value = -1806.6428
$ws.Range("K31").Value = 2101.6428
$ws.Range("M34").Value = -1899.6428
$ws.Range("K34").Value = 2101.6428
$ws.Range("I34").Value = 2101.6428
$ws.Range("H34").Value = 5116.8423
$ws.Range("H134").Value = 1786.3125
$ws.Range("I134").Value = 1649.0667
$ws.Range("K134").Value = 4947.2001
$ws.Range("M134").Value = -2412.2001

# --- Sheet: CUL (7 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M132").Value = -8987.427799999999
$ws.Range("I132").Value = 1279.7142
$ws.Range("K132").Value = 11517.4278
$ws.Range("N132").Value = -31847.9987
$ws.Range("L132").Value = 26787.9987
$ws.Range("H132").Value = 2234.125
$ws.Range("J132").Value = 2976.4443

# --- Sheet: GSM (49 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J18").Value = 4999
$ws.Range("N18").Value = -5585
$ws.Range("H18").Value = 4999
$ws.Range("L18").Value = 4999
$ws.Range("L43").Value = 6000
$ws.Range("J43").Value = 6000
$ws.Range("N43").Value = -6302
$ws.Range("H43").Value = 6600
$ws.Range("H46").Value = 25399.8
$ws.Range("K46").Value = 13999.667
$ws.Range("I46").Value = 13999.667
$ws.Range("M46").Value = -13843.667
$ws.Range("N52").Value = -224518
$ws.Range("J52").Value = 224000
$ws.Range("H52").Value = 224000
$ws.Range("L52").Value = 224000
$ws.Range("H57").Value = 10000000
$ws.Range("L57").Value = 10000000
$ws.Range("N57").Value = -10001640
$ws.Range("J57").Value = 10000000
$ws.Range("H80").Value = 8361.762000000001
$ws.Range("N80").Value = -6652.4165
$ws.Range("I80").Value = 13302.223
$ws.Range("L80").Value = 4656.4165
$ws.Range("J80").Value = 4656.4165
$ws.Range("M80").Value = -12304.223
$ws.Range("K80").Value = 13302.223
$ws.Range("L83").Value = 23282.0825
$ws.Range("J83").Value = 4656.4165
$ws.Range("I83").Value = 13302.223
$ws.Range("H83").Value = 8361.762000000001
$ws.Range("N83").Value = -33266.0825
$ws.Range("K83").Value = 66511.11500000001
$ws.Range("M83").Value = -61519.11500000001
$ws.Range("H122").Value = 8588.549000000001
$ws.Range("K122").Value = 27133.5
$ws.Range("M122").Value = -24683.5
$ws.Range("I122").Value = 9044.5
$ws.Range("H126").Value = 3377.6667
$ws.Range("I126").Value = 2659.4
$ws.Range("M126").Value = -5508.200000000001
$ws.Range("K126").Value = 7978.200000000001
$ws.Range("M132").Value = -6125.136200000001
$ws.Range("I132").Value = 2885.0454
$ws.Range("K132").Value = 8655.136200000001
$ws.Range("N132").Value = -14013.5
$ws.Range("L132").Value = 8953.5
$ws.Range("H132").Value = 2893.3333
$ws.Range("J132").Value = 2984.5

# --- Sheet: LTW (37 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -4433.6665
$ws.Range("I7").Value = 4545.6665
$ws.Range("H7").Value = 3838.7778
$ws.Range("K7").Value = 4545.6665
$ws.Range("J22").Value = 2026.9
$ws.Range("N22").Value = -2616.9
$ws.Range("K22").Value = 2153.8
$ws.Range("M22").Value = -1858.8
$ws.Range("L22").Value = 2026.9
$ws.Range("I22").Value = 2153.8
$ws.Range("H22").Value = 2069.2
$ws.Range("J27").Value = 2026.9
$ws.Range("M27").Value = -2046.8
$ws.Range("N27").Value = -2240.9
$ws.Range("L27").Value = 2026.9
$ws.Range("I27").Value = 2153.8
$ws.Range("K27").Value = 2153.8
$ws.Range("H27").Value = 2069.2
$ws.Range("H46").Value = 990.3333
$ws.Range("K46").Value = 974.5
$ws.Range("I46").Value = 974.5
$ws.Range("M46").Value = -786.5
$ws.Range("H126").Value = 3838.7778
$ws.Range("I126").Value = 4545.6665
$ws.Range("M126").Value = -11166.9995
$ws.Range("K126").Value = 13636.9995
$ws.Range("M132").Value = -6559.8287
$ws.Range("I132").Value = 3029.9429
$ws.Range("K132").Value = 9089.8287
$ws.Range("N132").Value = -19272.5
$ws.Range("L132").Value = 14212.5
$ws.Range("H132").Value = 3347.628
$ws.Range("J132").Value = 4737.5
$ws.Range("K136").Value = 9535.625100000001
$ws.Range("I136").Value = 3178.5417
$ws.Range("H136").Value = 3242.2856
$ws.Range("M136").Value = -6985.625100000001

# --- Sheet: WVR (11 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L122").Value = 25250.001
$ws.Range("N122").Value = -30150.001
$ws.Range("H122").Value = 3705.4583
$ws.Range("J122").Value = 8416.666999999999
$ws.Range("M132").Value = -4795.750100000001
$ws.Range("I132").Value = 2441.9167
$ws.Range("K132").Value = 7325.750100000001
$ws.Range("N132").Value = -11538.5
$ws.Range("L132").Value = 6478.5
$ws.Range("H132").Value = 2401.5715
$ws.Range("J132").Value = 2159.5
